$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.709.21'
$ws.Range("E2").Value = '  -4.09%  '

$ws.Range("D3").Value = '3.458.51'
$ws.Range("E3").Value = '  -4.26%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.28%  '

$ws.Range("D7").Value = '3.457.24'
$ws.Range("E7").Value = '  -4.24%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("E9").Value = '  -2.54%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.142'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.47'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.422'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.51%  '

$ws.Range("E13").Value = '  -7.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.69%  '

$ws.Range("D15").Value = '4.042.32'
$ws.Range("E15").Value = '  -4.26%  '

$ws.Range("D16").Value = '3.465.02'
$ws.Range("E16").Value = '  -3.97%  '

$ws.Range("D17").Value = '66.801.20'
$ws.Range("E17").Value = '  -3.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '438.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.15%  '

$ws.Range("E23").Value = '  -6.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.82%  '

$ws.Range("E25").Value = '  +0.04%  '

$ws.Range("D26").Value = '3.595.77'
$ws.Range("E26").Value = '  -4.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000120'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.91%  '

$ws.Range("E32").Value = '  -0.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.166'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.97%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.54%  '

$ws.Range("D36").Value = '3.452.21'
$ws.Range("E36").Value = '  -4.34%  '

$ws.Range("E37").Value = '  -7.86%  '

$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.16%  '

$ws.Range("E41").Value = '  -3.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.17'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0885'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.881'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.63%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.00%  '

$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.27%  '

$ws.Range("E48").Value = '  -11.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -11.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.986'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.67%  '
